$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.106.16'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.24%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.902.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.83%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5064'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3916'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.28%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09598'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.54%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.133'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.23%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.75'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.64%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.382'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.34%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.915.82'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.98%  '

$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.04%  '

$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.008'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.37%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.296'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001117'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.24%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06614'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.005'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.23%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.211'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.122.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.315'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.71%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.644'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.08%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.140.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.91%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.90%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '157.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.66%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.083'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.32%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1062'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.614'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.06%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.628'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.612'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.40%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06615'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.58%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02414'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.96%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.233'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.45%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2174'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.39%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.282'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.19%  '

$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.983'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.68%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6315'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.46%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.39'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.76%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.004'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.14%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.58%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5980'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.737'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.77%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.281'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.88%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.022'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.62%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.58%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.182'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.48%  '

Write-Host "Updated cryptos list"